# Add the new "2020" column (J) to the freshwater-withdrawal table, mirroring
# the formatting already used by the adjacent "2019" column (I).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column I's formatting (number formats, fonts, borders, alignment)
# into column J before writing the new figures into it.
$ws.Range("I4:I18").Copy()
$ws.Range("J4:J18").PasteSpecial(-4122)  # xlPasteFormats

# Year header
$ws.Range("J4").Value = 2020

# "Всего" (Total) block
$ws.Range("J5").Value = 8017.9
$ws.Range("J6").Value = ""

# "по видам источников" (by type of source) block
$ws.Range("J7").Formula = "=J5-J8"
$ws.Range("J8").Value = 249.8

# "по терретории" (by territory) block
$ws.Range("J9").Value = ""
$ws.Range("J10").Value = 757.6
$ws.Range("J11").Value = 984.4
$ws.Range("J12").Value = 646.2
$ws.Range("J13").Value = 667.6
$ws.Range("J14").Value = 1147
$ws.Range("J15").Value = 961.1
$ws.Range("J16").Value = 2664.5
$ws.Range("J17").Value = 132.5
$ws.Range("J18").Value = 57

# Match the post-edit selection left behind in the source workbook.
$null = $ws.Range("J19").Select()
